$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.349.50"
$ws.Range("D3").Value = "1.564.99"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.94"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.42"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.65"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "1.787.23"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "1.558.30"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "28.345.11"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.03"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.59"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.94"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.47"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "1.387.96"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.64"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0163"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.521"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.95"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.29"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.919"
$ws.Range("E48").Value = "  -5.96%  "
$ws.Range("D49").Value = "1.700.02"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  -2.05%  "
